# Applies the "Copy in EU-2024-develop branch" edit to the
# "Cap Ret per Unit Net Loss" (CRpUNL) workbook:
#   - About sheet: tweak the "reliability exemption" explanation text and
#     append two new paragraph lines about biomass/CHP co-location.
#   - CRpUNL sheet: change the header from a $/MWh-based fraction to a flat
#     "MW retired" @ "Unit: MW/($/MW)", and rewrite the retirement values -
#     plants that get an economic-retirement fraction now show 0.03, while
#     plants exempted from economic retirement (incl. now biomass) show 0
#     with the existing gray "exempt" shading.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "About"
# ---------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

$about.Range("A10").Value = "These includes: natural gas steam turbines and petroleum plants. For these plant types we set the "

# New explanatory lines about biomass/CHP, appended after the existing text.
$about.Range("A16").Value = "Likewise, biomass plants are often colocated with cheap supply and part of integrated"
$about.Range("A17").Value = "CHP or industrial systems, and we therefore do not subject them to economic retirement."

$about.Range("A18").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "CRpUNL"
# ---------------------------------------------------------------------
$crpunl = $wb.Worksheets.Item("CRpUNL")

$crpunl.Range("A1").Value = "Unit: MW/(`$/MW)"
$crpunl.Range("B1").Value = "MW retired"

# Plant types that keep getting economically retired: flat 0.03 MW/($/MW).
$economicRetirementCells = @("B2","B3","B4","B5","B7","B8","B13","B14","B15")
foreach ($cell in $economicRetirementCells) {
    $crpunl.Range($cell).Value = 0.03
}

# Plant types newly/still exempted from economic retirement: value 0, with
# the gray "exempt" shading copied over from a cell that already has it
# (B6 / hydro) so the style is reused rather than duplicated.
$crpunl.Range("B6").Copy() | Out-Null
$exemptCells = @("B9","B10","B11","B12","B18")
foreach ($cell in $exemptCells) {
    $crpunl.Range($cell).PasteSpecial(-4122) | Out-Null
    $crpunl.Range($cell).Value = 0
}

# Plant types that were already exempt (styled) but now get the 0.03 value.
$retirementFractionCells = @("B19","B20","B21","B22","B23","B24","B25")
foreach ($cell in $retirementFractionCells) {
    $crpunl.Range($cell).Value = 0.03
}

$crpunl.Range("D14").Select() | Out-Null
$about.Activate() | Out-Null
